$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9965318441390991
$ws.Range("B1").Value = 1.260233640670776
$ws.Range("C1").Value = 1.828446507453918
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 2.130897998809814
